$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Commit message: "Restored from revision ... SAVE" — the underlying content
# change is cell C10 on the Rules sheet, which reverts from 18 to 1.
$ws.Range("C10").Value = 1
